# Update TPM-derived edge metrics for the C1qa-Cd93 LR-pair sheet
# (Sending/Ligand/Receptor/Target cluster labels in columns A-D are unchanged;
#  only the computed expression/specificity values in columns E:T are refreshed
#  to reflect the new TPM input.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,16
$row2[0,0] = [double]"1"
$row2[0,1] = [double]"0.3333333333333333"
$row2[0,2] = [double]"0.1161373333333333"
$row2[0,3] = [double]"0.348412"
$row2[0,4] = [double]"0.0005423317774654872"
$row2[0,5] = [double]"0.0005423317774654872"
$row2[0,6] = [double]"3"
$row2[0,7] = [double]"1"
$row2[0,8] = [double]"209.26237"
$row2[0,9] = [double]"627.78711"
$row2[0,10] = [double]"0.8127157202241573"
$row2[0,11] = [double]"0.8127157202241573"
$row2[0,12] = [double]"24.30317361881333"
$row2[0,13] = [double]"218.72856256932"
$row2[0,14] = [double]"0.0004407615611233109"
$row2[0,15] = [double]"0.0004407615611233109"
$ws.Range("E2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,16
$row3[0,0] = [double]"1"
$row3[0,1] = [double]"0.3333333333333333"
$row3[0,2] = [double]"0.1161373333333333"
$row3[0,3] = [double]"0.348412"
$row3[0,4] = [double]"0.0005423317774654872"
$row3[0,5] = [double]"0.0005423317774654872"
$row3[0,6] = [double]"3"
$row3[0,7] = [double]"1"
$row3[0,8] = [double]"0.9848756666666668"
$row3[0,9] = [double]"2.954627"
$row3[0,10] = [double]"0.003824977881910862"
$row3[0,11] = [double]"0.003824977881910862"
$row3[0,12] = [double]"0.1143808335915556"
$row3[0,13] = [double]"1.029427502324"
$row3[0,14] = [double]"2.074407053462893e-06"
$row3[0,15] = [double]"2.074407053462893e-06"
$ws.Range("E3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,16
$row4[0,0] = [double]"1"
$row4[0,1] = [double]"0.3333333333333333"
$row4[0,2] = [double]"0.1161373333333333"
$row4[0,3] = [double]"0.348412"
$row4[0,4] = [double]"0.0005423317774654872"
$row4[0,5] = [double]"0.0005423317774654872"
$row4[0,6] = [double]"3"
$row4[0,7] = [double]"1"
$row4[0,8] = [double]"1.763846666666667"
$row4[0,9] = [double]"5.291539999999999"
$row4[0,10] = [double]"0.006850280411451801"
$row4[0,11] = [double]"0.006850280411451801"
$row4[0,12] = [double]"0.2048484482755555"
$row4[0,13] = [double]"1.84363603448"
$row4[0,14] = [double]"3.715124751679664e-06"
$row4[0,15] = [double]"3.715124751679664e-06"
$ws.Range("E4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,16
$row5[0,0] = [double]"1"
$row5[0,1] = [double]"0.3333333333333333"
$row5[0,2] = [double]"0.1161373333333333"
$row5[0,3] = [double]"0.348412"
$row5[0,4] = [double]"0.0005423317774654872"
$row5[0,5] = [double]"0.0005423317774654872"
$row5[0,6] = [double]"3"
$row5[0,7] = [double]"1"
$row5[0,8] = [double]"45.474231"
$row5[0,9] = [double]"136.422693"
$row5[0,10] = [double]"0.1766090214824801"
$row5[0,11] = [double]"0.1766090214824801"
$row5[0,12] = [double]"5.281255923723999"
$row5[0,13] = [double]"47.53130331351599"
$row5[0,14] = [double]"9.578068453703385e-05"
$row5[0,15] = [double]"9.578068453703385e-05"
$ws.Range("E5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,16
$row6[0,0] = [double]"3"
$row6[0,1] = [double]"1"
$row6[0,2] = [double]"214.0283"
$row6[0,3] = [double]"642.0849000000001"
$row6[0,4] = [double]"0.9994576682225345"
$row6[0,5] = [double]"0.9994576682225345"
$row6[0,6] = [double]"3"
$row6[0,7] = [double]"1"
$row6[0,8] = [double]"209.26237"
$row6[0,9] = [double]"627.78711"
$row6[0,10] = [double]"0.8127157202241573"
$row6[0,11] = [double]"0.8127157202241573"
$row6[0,12] = [double]"44788.06930507101"
$row6[0,13] = [double]"403092.6237456391"
$row6[0,14] = [double]"0.8122749586630339"
$row6[0,15] = [double]"0.8122749586630339"
$ws.Range("E6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,16
$row7[0,0] = [double]"3"
$row7[0,1] = [double]"1"
$row7[0,2] = [double]"214.0283"
$row7[0,3] = [double]"642.0849000000001"
$row7[0,4] = [double]"0.9994576682225345"
$row7[0,5] = [double]"0.9994576682225345"
$row7[0,6] = [double]"3"
$row7[0,7] = [double]"1"
$row7[0,8] = [double]"0.9848756666666668"
$row7[0,9] = [double]"2.954627"
$row7[0,10] = [double]"0.003824977881910862"
$row7[0,11] = [double]"0.003824977881910862"
$row7[0,12] = [double]"210.7912646480334"
$row7[0,13] = [double]"1897.1213818323"
$row7[0,14] = [double]"0.003822903474857399"
$row7[0,15] = [double]"0.003822903474857399"
$ws.Range("E7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,16
$row8[0,0] = [double]"3"
$row8[0,1] = [double]"1"
$row8[0,2] = [double]"214.0283"
$row8[0,3] = [double]"642.0849000000001"
$row8[0,4] = [double]"0.9994576682225345"
$row8[0,5] = [double]"0.9994576682225345"
$row8[0,6] = [double]"3"
$row8[0,7] = [double]"1"
$row8[0,8] = [double]"1.763846666666667"
$row8[0,9] = [double]"5.291539999999999"
$row8[0,10] = [double]"0.006850280411451801"
$row8[0,11] = [double]"0.006850280411451801"
$row8[0,12] = [double]"377.5131035273334"
$row8[0,13] = [double]"3397.617931746"
$row8[0,14] = [double]"0.006846565286700121"
$row8[0,15] = [double]"0.006846565286700121"
$ws.Range("E8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,16
$row9[0,0] = [double]"3"
$row9[0,1] = [double]"1"
$row9[0,2] = [double]"214.0283"
$row9[0,3] = [double]"642.0849000000001"
$row9[0,4] = [double]"0.9994576682225345"
$row9[0,5] = [double]"0.9994576682225345"
$row9[0,6] = [double]"3"
$row9[0,7] = [double]"1"
$row9[0,8] = [double]"45.474231"
$row9[0,9] = [double]"136.422693"
$row9[0,10] = [double]"0.1766090214824801"
$row9[0,11] = [double]"0.1766090214824801"
$row9[0,12] = [double]"9732.7723547373"
$row9[0,13] = [double]"87594.9511926357"
$row9[0,14] = [double]"0.176513240797943"
$row9[0,15] = [double]"0.176513240797943"
$ws.Range("E9:T9").Value = $row9

Write-Host "Updated TPM values for C1qa-Cd93 rows 2-9"
